$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (VillagerA) - fix prefab path + tune stats
$ws.Range("C4").Value = "Assets/Prefabs/NPC/VillagerA.prefab"
$ws.Range("E4").Value = 10.0
$ws.Range("G4").Value = 0.7

# Row 5 now represents the Warrior (replacing old VillagerB data)
$ws.Range("A5").Value = "Warrior"
$ws.Range("B5").Value = 20.0
$ws.Range("C5").Value = "Assets/Prefabs/NPC/Warrior.prefab"
$ws.Range("D5").Value = "WarriorStats"
$ws.Range("E5").Value = 200.0
$ws.Range("F5").Value = 100.0
$ws.Range("I5").Value = "Warrior"

# Remove the now-unused blank row 11 (keeps row 12 labeled as 12)
$ws.Range("A11:I11").Clear()
